# Apply updated "dSF" (column F) values coming from a data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 5
    4  = 1
    5  = -1
    6  = -2
    7  = 5
    8  = 2
    9  = 4
    10 = 1
    11 = 7
    13 = -5
    14 = -1
    15 = 5
    16 = 2
    17 = -2
    18 = 3
    19 = -1
    21 = -1
    22 = 1
    23 = 4
    24 = -3
    25 = -5
    27 = 1
    28 = -5
    29 = 5
    30 = -4
    31 = -3
    32 = -6
    33 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
